$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$ws.Cells.Item(2, 7).Value = "mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Cells.Item(3, 7).Value = "nourhan.mahmoud@med.asu.edu.eg , shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Cells.Item(9, 7).Value = "norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Cells.Item(14, 7).Value = "Alshimaa_khaled@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Cells.Item(15, 7).Value = "mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg"
$ws.Cells.Item(27, 7).Value = "randa.rabea@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg, System"
$ws.Cells.Item(29, 7).Value = "nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg, marina_atef@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
$ws.Cells.Item(33, 7).Value = "mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Cells.Item(34, 7).Value = "nourhan.mahmoud@med.asu.edu.eg , shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Cells.Item(40, 7).Value = "norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Cells.Item(45, 7).Value = "Alshimaa_khaled@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Cells.Item(46, 7).Value = "mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"
$ws.Cells.Item(58, 7).Value = "randa.rabea@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg, System"
$ws.Cells.Item(60, 7).Value = "nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg, marina_atef@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
$ws.Cells.Item(64, 7).Value = "Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Cells.Item(65, 7).Value = "Amira.Sobhy@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Cells.Item(71, 7).Value = "Omnia.Mohammed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Cells.Item(76, 7).Value = "mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg"
$ws.Cells.Item(77, 7).Value = "mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"
$ws.Cells.Item(91, 7).Value = "nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg, marina_atef@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
$ws.Cells.Item(95, 7).Value = "Veronia.rafat@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Cells.Item(96, 7).Value = "alshimaa.atef@med.asu.edu.egm, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg "
$ws.Cells.Item(102, 7).Value = "norhan.mohamed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Cells.Item(107, 7).Value = "Alshimaa_khaled@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Cells.Item(108, 7).Value = "mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg"
$ws.Cells.Item(120, 7).Value = "dr_ahmed_mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Cells.Item(122, 7).Value = "naema.gomaa@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Cells.Item(126, 7).Value = "Veronia.rafat@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Cells.Item(127, 7).Value = "alshimaa.atef@med.asu.edu.egm, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg "
$ws.Cells.Item(133, 7).Value = "norhan.mohamed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Cells.Item(138, 7).Value = "mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg"
$ws.Cells.Item(139, 7).Value = "mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Cells.Item(140, 7).Value = "mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"
$ws.Cells.Item(151, 7).Value = "Amr-Saeed@med.asu.edu.eg`n, takyas@med.asu.edu.eg, enas.omran@med.asu.edu.eg, user@user.com"
$ws.Cells.Item(153, 7).Value = "naema.gomaa@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Cells.Item(157, 7).Value = "majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Cells.Item(158, 7).Value = "Amira.Sobhy@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Cells.Item(164, 7).Value = "Omnia.Mohammed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Cells.Item(169, 7).Value = "Alshimaa_khaled@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Cells.Item(170, 7).Value = "mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Cells.Item(182, 7).Value = "dr_ahmed_mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Cells.Item(184, 7).Value = "ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
